$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (About!C1) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45379

# --- "FPIEBP" sheet: update "hard coal" balancing priorities (row 3) ---
$fpiebp = $wb.Worksheets.Item("FPIEBP")
$fpiebp.Range("B3").Value = 1
$fpiebp.Range("C3").Value = 3
$fpiebp.Range("D3").Value = 2

# --- Restore the active selection on the FPIEBP sheet ---
$fpiebp.Activate()
$fpiebp.Range("E3").Select()
